$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 15 de Junio de 2020 a las 19:20"

# Swap country labels where new data caused a reordering/overwrite of adjacent rows
$ws.Range("A14").Value = "Chile"
$ws.Range("A15").Value = "Turquia"
$ws.Range("A67").Value = "Marruecos"
$ws.Range("A68").Value = "Honduras"
$ws.Range("A179").Value = "Barbados"
$ws.Range("A180").Value = "Eritrea"
$ws.Range("A206").Value = "Islas Malvinas"
$ws.Range("A207").Value = "Groenlandia"
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("A209").Value = "Islas Turcas y Caicos"

# Apply updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B4").Value = 2170722
$ws.Range("C4").Value = 8494
$ws.Range("D4").Value = 873753
$ws.Range("E4").Value = 1179010
$ws.Range("G4").Value = 101
$ws.Range("H4").Value = 117959
$ws.Range("B5").Value = 873963
$ws.Range("C5").Value = 6081
$ws.Range("E5").Value = 376910
$ws.Range("G5").Value = 96
$ws.Range("H5").Value = 43485
$ws.Range("B7").Value = 342291
$ws.Range("C7").Value = 9508
$ws.Range("D7").Value = 179321
$ws.Range("E7").Value = 153086
$ws.Range("G7").Value = 364
$ws.Range("H7").Value = 9884
$ws.Range("B10").Value = 237290
$ws.Range("C10").Value = 301
$ws.Range("D10").Value = 177010
$ws.Range("E10").Value = 25909
$ws.Range("G10").Value = 26
$ws.Range("H10").Value = 34371
$ws.Range("B13").Value = 187854
$ws.Range("C13").Value = 183
$ws.Range("E13").Value = 6379
$ws.Range("B14").Value = 179436
$ws.Range("C14").Value = 5143
$ws.Range("D14").Value = 143704
$ws.Range("E14").Value = 32370
$ws.Range("G14").Value = 39
$ws.Range("H14").Value = 3362
$ws.Range("B15").Value = 178239
$ws.Range("D15").Value = 151417
$ws.Range("E15").Value = 22015
$ws.Range("H15").Value = 4807
$ws.Range("D39").Value = 28900
$ws.Range("E39").Value = 292
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 1939
$ws.Range("E42").Value = 19885
$ws.Range("G42").Value = 7
$ws.Range("H42").Value = 478
$ws.Range("B43").Value = 25321
$ws.Range("C43").Value = 18
$ws.Range("E43").Value = 917
$ws.Range("B62").Value = 11031
$ws.Range("C62").Value = 112
$ws.Range("D62").Value = 7735
$ws.Range("E62").Value = 2519
$ws.Range("G62").Value = 10
$ws.Range("H62").Value = 777
$ws.Range("B67").Value = 8885
$ws.Range("C67").Value = 92
$ws.Range("D67").Value = 7828
$ws.Range("E67").Value = 845
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 212
$ws.Range("B68").Value = 8858
$ws.Range("C68").Value = 403
$ws.Range("D68").Value = 967
$ws.Range("E68").Value = 7579
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 312
$ws.Range("D76").Value = 4019
$ws.Range("E76").Value = 1116
$ws.Range("B81").Value = 4501
$ws.Range("C81").Value = 36
$ws.Range("D81").Value = 3183
$ws.Range("E81").Value = 1275
$ws.Range("B85").Value = 4072
$ws.Range("C85").Value = 2
$ws.Range("D85").Value = 3931
$ws.Range("B95").Value = 2642
$ws.Range("C95").Value = 24
$ws.Range("D95").Value = 622
$ws.Range("E95").Value = 1932
$ws.Range("B143").Value = 609
$ws.Range("C143").Value = 26
$ws.Range("D143").Value = 157
$ws.Range("E143").Value = 449
$ws.Range("B179").Value = 97
$ws.Range("C179").Value = 1
$ws.Range("D179").Value = 83
$ws.Range("E179").Value = 7
$ws.Range("H179").Value = 7
$ws.Range("D180").Value = 39
$ws.Range("E180").Value = 57
$ws.Range("H180").Value = 0
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1
